# Weekly fruit/vegetable data update: a new record for this market/product
# was collected, so it is inserted at the top of the data block (row 62),
# pushing all the existing records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 62 (shifts old rows 62:96 down to 63:97,
# and extends the used range to A1:R97).
$ws.Rows(62).Insert()

# Populate the newly inserted row 62 with this week's record.
$ws.Cells.Item(62, 1).Value = 8
$ws.Cells.Item(62, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(62, 3).Value = "Coquimbo"
$ws.Cells.Item(62, 4).Value = 44518
$ws.Cells.Item(62, 5).Value = 4
$ws.Cells.Item(62, 6).Value = 100112044
$ws.Cells.Item(62, 7).Value = "Perejil"
$ws.Cells.Item(62, 8).Value = "Sin especificar"
$ws.Cells.Item(62, 9).Value = "Primera"
$ws.Cells.Item(62, 10).Value = 3320
$ws.Cells.Item(62, 11).Value = 1300
$ws.Cells.Item(62, 12).Value = 1500
$ws.Cells.Item(62, 13).Value = 1400
$ws.Cells.Item(62, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(62, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(62, 16).Value = 933
$ws.Cells.Item(62, 17).Value = 1.5
$ws.Cells.Item(62, 18).Value = "Hortaliza"
